# Edit: replace the verbose "https://github.com/.../CAP-5610_Machine-Learning.git"
# run sequence on the GitHub slide with a single hyperlinked "GitHub Link" run,
# per commit message "Changed the link in the presentation".

$p = $ppt.ActivePresentation

# Locate the slide / shape that holds "My GitHub link:" and the raw URL text.
$targetSlide = $null
$targetShape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shp = $slide.Shapes.Item($shi)
        if ($shp.HasTextFrame) {
            $t = $shp.TextFrame.TextRange.Text
            if ($t -and $t.IndexOf("CAP-5610_Machine-Learning.git") -ge 0) {
                $targetSlide = $slide
                $targetShape = $shp
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# Find the "https://...git " run of text (it is immediately followed by the
# paragraph mark of its own paragraph, and that paragraph is followed by one
# empty paragraph which collapses away once the text run + break are removed).
$full = $tr.Text
$urlText = "https://github.com/monicabernard/CAP-5610_Machine-Learning.git "
$idx = $full.IndexOf($urlText)

# Select the URL text plus the paragraph break right after it so replacing it
# merges the now-empty trailing paragraph into this one (matching the diff,
# which drops one paragraph).
$selLen = $urlText.Length + 1
$urlRange = $tr.Characters($idx + 1, $selLen)
$urlRange.Text = "GitHub Link"

# Re-acquire the newly written "GitHub Link" run and point it at the GitHub
# repo as an external hyperlink.
$full2 = $tr.Text
$linkText = "GitHub Link"
$idx2 = $full2.IndexOf($linkText)
$linkRange = $tr.Characters($idx2 + 1, $linkText.Length)

$actionSetting = $linkRange.ActionSettings.Item(1)
$actionSetting.Hyperlink.Address = "https://github.com/monicabernard/CAP-5610_Machine-Learning.git"
